$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns B, C, D are treated as text so numeric-looking strings
# (e.g. "211.80", "1.00", "0.102") are not auto-converted to numbers by Excel.
$ws.Range("B2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.676.37'
$ws.Range("E2").Value = '  -0.14%  '

$ws.Range("D3").Value = '1.599.40'
$ws.Range("E3").Value = '  +0.05%  '

$ws.Range("E4").Value = '  +0.31%  '

$ws.Range("D5").Value = '211.80'
$ws.Range("E5").Value = '  -0.03%  '

$ws.Range("E6").Value = '  +0.55%  '

$ws.Range("E7").Value = '  +0.28%  '

$ws.Range("D8").Value = '0.0619'
$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("D9").Value = '0.247'
$ws.Range("E9").Value = '  +0.25%  '

$ws.Range("E10").Value = '  -0.81%  '

$ws.Range("E11").Value = '  +0.41%  '

$ws.Range("D12").Value = '1.823.98'
$ws.Range("E12").Value = '  +0.01%  '

$ws.Range("D13").Value = '1.590.99'
$ws.Range("E13").Value = '  -0.23%  '

$ws.Range("E14").Value = '  +0.07%  '

$ws.Range("E15").Value = '  +0.07%  '

$ws.Range("D16").Value = '65.11'
$ws.Range("E16").Value = '  +0.05%  '

$ws.Range("D17").Value = '26.675.93'
$ws.Range("E17").Value = '  -0.04%  '

$ws.Range("D18").Value = '0.0₃0737'
$ws.Range("E18").Value = '  +1.15%  '

$ws.Range("E19").Value = '  +0.37%  '

$ws.Range("D20").Value = '208.66'
$ws.Range("E20").Value = '  -0.63%  '

$ws.Range("D21").Value = '7.05'
$ws.Range("E21").Value = '  +5.07%  '

$ws.Range("E22").Value = '  +0.63%  '

$ws.Range("E23").Value = '  +0.59%  '

$ws.Range("D24").Value = '8.98'
$ws.Range("E24").Value = '  +0.74%  '

$ws.Range("D25").Value = '145.19'
$ws.Range("E25").Value = '  -1.17%  '

$ws.Range("E26").Value = '  +0.22%  '

$ws.Range("E27").Value = '  -0.63%  '

$ws.Range("E28").Value = '  -0.36%  '

$ws.Range("D29").Value = '15.31'
$ws.Range("E29").Value = '  -0.07%  '

$ws.Range("D30").Value = '0.0514'
$ws.Range("E30").Value = '  +2.15%  '

$ws.Range("E31").Value = '  +0.11%  '

$ws.Range("E32").Value = '  +0.42%  '

$ws.Range("D33").Value = '2.94'
$ws.Range("E33").Value = '  +1.28%  '

$ws.Range("D34").Value = '1.276.81'
$ws.Range("E34").Value = '  -1.66%  '

$ws.Range("E35").Value = '  -7.68%  '

$ws.Range("E36").Value = '  +0.35%  '

$ws.Range("E37").Value = '  +0.94%  '

$ws.Range("E38").Value = '  -0.88%  '

$ws.Range("E39").Value = '  -0.84%  '

$ws.Range("E40").Value = '  +18.54%  '

$ws.Range("E41").Value = '  +2.45%  '

$ws.Range("E42").Value = '  +0.25%  '

$ws.Range("E43").Value = '  -0.76%  '

$ws.Range("E44").Value = '  +0.14%  '

$ws.Range("D45").Value = '1.736.05'
$ws.Range("E45").Value = '  -0.01%  '

$ws.Range("D46").Value = '91.19'
$ws.Range("E46").Value = '  +1.33%  '

$ws.Range("E47").Value = '  -2.41%  '

$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").Value = '0.102'
$ws.Range("E48").Value = '  +3.26%  '

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.0508'
$ws.Range("E49").Value = '  +0.61%  '

$ws.Range("B50").Value = 'USDD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D50").Value = '1.00'
$ws.Range("E50").Value = '  +0.04%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '7.38'
$ws.Range("E51").Value = '  -1.75%  '
